$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9989239189189486
$ws.Range("D2").Value = 0.3248881816885114

$ws.Range("C3").Value = -0.1431471177330473
$ws.Range("D3").Value = 0.887019006442983

$ws.Range("C4").Value = 0.2447833943557833
$ws.Range("D4").Value = 0.8080955039427922

$ws.Range("C5").Value = -0.2332024512668203
$ws.Range("D5").Value = 0.8170025101505642

$ws.Range("C6").Value = -0.7395626055895445
$ws.Range("D6").Value = 0.4646413962346965

$ws.Range("C7").Value = -0.2590954419234484
$ws.Range("D7").Value = 0.7971238990803957

$ws.Range("C8").Value = -0.6819652277681563
$ws.Range("D8").Value = 0.4998806459126177

$ws.Range("C9").Value = 0.5251714458173908
$ws.Range("D9").Value = 0.6028721944047715

$ws.Range("C10").Value = -0.1296090855199939
$ws.Range("D10").Value = 0.8976394587219785

$ws.Range("C11").Value = -0.3041759945811227
$ws.Range("D11").Value = 0.7628477415070645
